$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value  = 173
$ws.Range("I2").Value  = 532
$ws.Range("J2").Value  = 2121
$ws.Range("K2").Value  = 11
$ws.Range("L2").Value  = 579
$ws.Range("M2").Value  = 41
$ws.Range("N2").Value  = 381
$ws.Range("P2").Value  = 10
$ws.Range("Q2").Value  = 6
$ws.Range("R2").Value  = 26
$ws.Range("S2").Value  = 249
$ws.Range("T2").Value  = 351
$ws.Range("U2").Value  = 25
$ws.Range("V2").Value  = 3287
$ws.Range("W2").Value  = 2
$ws.Range("X2").Value  = 3285
$ws.Range("Z2").Value  = 56
$ws.Range("AA2").Value = 18

$wb.Save()
